$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.418.82"
$ws.Range("E2").Value = "  -1.84%  "

$ws.Range("D3").Value = "'3.870.50"
$ws.Range("E3").Value = "  -2.71%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'597.99"
$ws.Range("E5").Value = "  +1.63%  "

$ws.Range("D6").Value = "'167.69"
$ws.Range("E6").Value = "  +6.17%  "

$ws.Range("D7").Value = "'0.675"
$ws.Range("E7").Value = "  -0.57%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").Value = "'0.754"
$ws.Range("E9").Value = "  +0.75%  "

$ws.Range("D10").Value = "'0.177"
$ws.Range("E10").Value = "  +5.75%  "

$ws.Range("D11").Value = "'53.38"
$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("D12").Value = "'0.0000321"
$ws.Range("E12").Value = "  +1.32%  "

$ws.Range("D13").Value = "'11.19"
$ws.Range("E13").Value = "  +3.18%  "

$ws.Range("D14").Value = "'4.512.20"
$ws.Range("E14").Value = "  -2.50%  "

$ws.Range("D15").Value = "'3.902.62"
$ws.Range("E15").Value = "  -1.87%  "

$ws.Range("D16").Value = "'20.80"
$ws.Range("E16").Value = "  +1.89%  "

$ws.Range("D17").Value = "'13.87"
$ws.Range("E17").Value = "  -1.00%  "

$ws.Range("D18").Value = "'1.21"
$ws.Range("E18").Value = "  -5.12%  "

$ws.Range("E19").Value = "  -1.90%  "

$ws.Range("D20").Value = "'71.471.59"
$ws.Range("E20").Value = "  -1.55%  "

$ws.Range("D21").Value = "'434.48"
$ws.Range("E21").Value = "  +0.85%  "

$ws.Range("D22").Value = "'4.75"
$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("D23").Value = "'94.47"
$ws.Range("E23").Value = "  -1.53%  "

$ws.Range("D24").Value = "'3.28"
$ws.Range("E24").Value = "  -3.86%  "

$ws.Range("D25").Value = "'13.81"
$ws.Range("E25").Value = "  -3.47%  "

$ws.Range("D26").Value = "'4.06"
$ws.Range("E26").Value = "  -7.78%  "

$ws.Range("D27").Value = "'10.90"
$ws.Range("E27").Value = "  -4.03%  "

$ws.Range("D28").Value = "'5.94"
$ws.Range("E28").Value = "  +0.19%  "

$ws.Range("D29").Value = "'10.18"
$ws.Range("E29").Value = "  -5.88%  "

$ws.Range("D30").Value = "'35.01"
$ws.Range("E30").Value = "  -3.75%  "

$ws.Range("D31").Value = "'7.81"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("D32").Value = "'50.51"
$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("D33").Value = "'13.55"
$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("D34").Value = "'0.125"
$ws.Range("E34").Value = "  -4.55%  "

$ws.Range("D35").Value = "'0.0₃0987"
$ws.Range("E35").Value = "  +15.26%  "

$ws.Range("D36").Value = "'68.93"
$ws.Range("E36").Value = "  +0.33%  "

$ws.Range("D37").Value = "'617.31"
$ws.Range("E37").Value = "  -9.16%  "

$ws.Range("D38").Value = "'0.419"
$ws.Range("E38").Value = "  -4.36%  "

$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "'3.27"
$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.142"
$ws.Range("E42").Value = "  -1.84%  "

$ws.Range("D43").Value = "'3.18"
$ws.Range("E43").Value = "  +32.49%  "

$ws.Range("D44").Value = "'0.0468"
$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("D45").Value = "'10.21"
$ws.Range("E45").Value = "  -6.83%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.144"
$ws.Range("E46").Value = "  -3.26%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.63"
$ws.Range("E47").Value = "  -2.36%  "

$ws.Range("D48").Value = "'3.38"
$ws.Range("E48").Value = "  -0.82%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'2.840.86"
$ws.Range("E49").Value = "  +2.70%  "

$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "'2.75"
$ws.Range("E50").Value = "  -18.93%  "

$ws.Range("D51").Value = "'0.000271"
$ws.Range("E51").Value = "  +0.48%  "
